# Auto-update draw results: append the new day's Pick 3 result as the
# next row right after the last existing data row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# Values such as "2025-09-28" and "250928" look numeric/date-like to
# Excel's auto-detection, so prefix them with a leading apostrophe to
# force plain text entry (keeps them identical in kind to every other
# row in the table, which stores Date/Phase as text, not numbers).
$ws.Range("A" + $newRow).Value = "'2025-09-28"
$ws.Range("B" + $newRow).Value = "Pick 3"
$ws.Range("C" + $newRow).Value = "'250928"
$ws.Range("D" + $newRow).Value = "6-9-2"
$ws.Range("E" + $newRow).Value = "2025-09-28T21:34:49.973+04:00"

# The apostrophe-prefix trick marks the cells with a "quote prefix"
# style so Excel remembers they were typed as text; reset the row's
# style back to Normal so it matches the unstyled cells around it.
$ws.Range("A" + $newRow + ":E" + $newRow).Style = "Normal"
